$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L10").Value = 92.53
$wsGrupo.Range("L17").Value = 21.02
$wsGrupo.Range("L30").Value = "4 de 28"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F10").Value = 92.53
$wsMensual.Range("F17").Value = 21.02
$wsMensual.Range("F30").Value = 1096.8

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# ColumnWidth property adds an internal padding of 5/6 (0.8333...) relative to the
# raw OOXML "width" attribute, so subtract it to land exactly on width=13 in the XML.
$wsCumpl.Columns.Item(4).ColumnWidth = 12.166666666666666
$wsCumpl.Range("D16").Value = 1069.07
$wsCumpl.Range("E16").Value = 17729.54
$wsCumpl.Range("F16").Value = 0.05686963025457733
$wsCumpl.Range("D19").Value = 1096.8
$wsCumpl.Range("E19").Value = 28440.99107555787
$wsCumpl.Range("F19").Value = 0.03713209282286471
